# Adafruit IO sync added one more reading to the feed log: append it as
# the new last row (row 82), extending the used range from A1:F81 to A1:F82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
# Value looks numeric ("25"); prefix with an apostrophe so Excel keeps it as
# text (matching every other Value/Lat/Lon/Elevation cell in this sheet),
# then drop the resulting quote-prefix formatting so no stray per-cell style
# is left behind.
$ws.Cells.Item($row, 3).Value = "'25"
$ws.Cells.Item($row, 3).ClearFormats()
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
